$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep column D as text so numeric-looking price strings (e.g. "1.00", "0.150")
# retain their exact literal formatting instead of being coerced to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '81.634.90'
$ws.Range('E2').Value = '  +2.26%  '
$ws.Range('D3').Value = '3.153.97'
$ws.Range('E3').Value = '  -1.56%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = '216.83'
$ws.Range('E5').Value = '  +3.96%  '
$ws.Range('D6').Value = '615.57'
$ws.Range('E6').Value = '  -4.14%  '
$ws.Range('D7').Value = '0.285'
$ws.Range('E7').Value = '  +15.05%  '
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').Value = '0.577'
$ws.Range('E9').Value = '  -5.29%  '
$ws.Range('D10').Value = '3.148.40'
$ws.Range('E10').Value = '  -1.71%  '
$ws.Range('D11').Value = '0.585'
$ws.Range('E11').Value = '  -5.36%  '
$ws.Range('D12').Value = '0.0000252'
$ws.Range('E12').Value = '  -7.44%  '
$ws.Range('E13').Value = '  -0.91%  '
$ws.Range('D14').Value = '5.28'
$ws.Range('E14').Value = '  -2.90%  '
$ws.Range('D15').Value = '3.736.14'
$ws.Range('E15').Value = '  -1.49%  '
$ws.Range('D16').Value = '31.83'
$ws.Range('E16').Value = '  -2.49%  '
$ws.Range('D17').Value = '81.779.19'
$ws.Range('E17').Value = '  +2.79%  '
$ws.Range('D18').Value = '3.163.14'
$ws.Range('E18').Value = '  -1.18%  '
$ws.Range('D19').Value = '3.20'
$ws.Range('E19').Value = '  +7.10%  '
$ws.Range('D20').Value = '13.86'
$ws.Range('E20').Value = '  -5.76%  '
$ws.Range('D21').Value = '434.07'
$ws.Range('E21').Value = '  -2.68%  '
$ws.Range('D22').Value = '8.81'
$ws.Range('E22').Value = '  -7.42%  '
$ws.Range('D23').Value = '5.10'
$ws.Range('E23').Value = '  -5.22%  '
$ws.Range('E24').Value = '  +4.23%  '
$ws.Range('D25').Value = '5.16'
$ws.Range('E25').Value = '  +5.88%  '
$ws.Range('D26').Value = '11.81'
$ws.Range('E26').Value = '  +7.28%  '
$ws.Range('B27').Value = 'WrappedeETH'
$ws.Range('C27').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D27').Value = '3.313.60'
$ws.Range('E27').Value = '  -1.63%  '
$ws.Range('B28').Value = 'Litecoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D28').Value = '76.24'
$ws.Range('E28').Value = '  -2.28%  '
$ws.Range('D29').Value = '0.998'
$ws.Range('E29').Value = '  +0.05%  '
$ws.Range('D30').Value = '0.0000119'
$ws.Range('E30').Value = '  -6.63%  '
$ws.Range('E31').Value = '  -0.02%  '
$ws.Range('D32').Value = '8.92'
$ws.Range('E32').Value = '  -3.61%  '
$ws.Range('D33').Value = '560.01'
$ws.Range('E33').Value = '  +1.50%  '
$ws.Range('D34').Value = '1.47'
$ws.Range('E34').Value = '  -5.45%  '
$ws.Range('D35').Value = '0.146'
$ws.Range('E35').Value = '  +18.41%  '
$ws.Range('D36').Value = '0.150'
$ws.Range('E36').Value = '  -5.33%  '
$ws.Range('D37').Value = '1.97'
$ws.Range('E37').Value = '  -4.33%  '
$ws.Range('D38').Value = '22.47'
$ws.Range('E38').Value = '  -3.81%  '
$ws.Range('B39').Value = 'FirstDigitalUSD'
$ws.Range('C39').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D39').Value = '0.999'
$ws.Range('E39').Value = '  -0.06%  '
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D40').Value = '6.14'
$ws.Range('E40').Value = '  +8.14%  '
$ws.Range('D41').Value = '0.402'
$ws.Range('E41').Value = '  -4.08%  '
$ws.Range('D42').Value = '20.82'
$ws.Range('E42').Value = '  +3.87%  '
$ws.Range('D43').Value = '1.99'
$ws.Range('E43').Value = '  +7.35%  '
$ws.Range('D44').Value = '2.98'
$ws.Range('E44').Value = '  +8.17%  '
$ws.Range('D45').Value = '158.65'
$ws.Range('E45').Value = '  -3.98%  '
$ws.Range('E46').Value = '  +0.05%  '
$ws.Range('D47').Value = '184.80'
$ws.Range('E47').Value = '  -5.48%  '
$ws.Range('D48').Value = '44.38'
$ws.Range('E48').Value = '  +0.75%  '
$ws.Range('D49').Value = '1.30'
$ws.Range('E49').Value = '  -3.35%  '
$ws.Range('B50').Value = 'Mantle'
$ws.Range('C50').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D50').Value = '0.758'
$ws.Range('E50').Value = '  -6.47%  '
$ws.Range('B51').Value = 'InjectiveProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D51').Value = '25.59'
$ws.Range('E51').Value = '  -2.67%  '
